# faturamento_diario.xlsx — "atualizei dados da bibi e add"
#
# Changes:
#  1) July (07/2025) day 2 total_venda (B3) corrected: 20464.65 -> 21194.65
#  2) A new row is inserted for July (07/2025) day 3, total_venda 22837.33,
#     right after the existing July day-2 row (old row 4), pushing the June
#     block (and everything after it) down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the existing July day-2 sale total.
$ws.Range("B3").Value = 21194.65

# Insert a new row at position 4 (shifts old row 4.. down by one) for the
# new July day-3 entry.
$ws.Rows.Item(4).Insert()

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 22837.33
$ws.Range("C4").Value = 7
$ws.Range("D4").Value = 2025
$ws.Range("E4").Value = "07/2025"
